# Update Customer_template.xlsx
#  1. Rename Sheet1 -> Template
#  2. Insert a new "Instructions" sheet right after Template
#  3. Expand the single credit_limits.* / sales_team.* child-table columns
#     on the Template sheet into 5 repeated column-groups each, and move
#     the remaining top-level fields after them.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Template"

# --- Rebuild header row on the Template sheet --------------------------

# First, stamp the header style (bold font / border / center-top align,
# i.e. the same style already used by the existing header cells) across
# the full new header range so every header cell - old and new - shares
# that formatting.
$ws1.Range("A1").Copy()
$ws1.Range("A1:AE1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A1").Value  = "customer_name [Data]"
$ws1.Range("B1").Value  = "customer_type [Select]"
$ws1.Range("C1").Value  = "is_internal_customer [Check]"
$ws1.Range("D1").Value  = "website [Data]"
$ws1.Range("E1").Value  = "customer_details [Text]"
$ws1.Range("F1").Value  = "tax_id [Data]"

$ws1.Range("G1").Value  = "credit_limits.1.credit_limit [Currency]"
$ws1.Range("H1").Value  = "credit_limits.1.bypass_credit_limit_check [Check]"
$ws1.Range("I1").Value  = "credit_limits.2.credit_limit [Currency]"
$ws1.Range("J1").Value  = "credit_limits.2.bypass_credit_limit_check [Check]"
$ws1.Range("K1").Value  = "credit_limits.3.credit_limit [Currency]"
$ws1.Range("L1").Value  = "credit_limits.3.bypass_credit_limit_check [Check]"
$ws1.Range("M1").Value  = "credit_limits.4.credit_limit [Currency]"
$ws1.Range("N1").Value  = "credit_limits.4.bypass_credit_limit_check [Check]"
$ws1.Range("O1").Value  = "credit_limits.5.credit_limit [Currency]"
$ws1.Range("P1").Value  = "credit_limits.5.bypass_credit_limit_check [Check]"

$ws1.Range("Q1").Value  = "sales_team.1.allocated_percentage [Float]"
$ws1.Range("R1").Value  = "sales_team.1.incentives [Currency]"
$ws1.Range("S1").Value  = "sales_team.2.allocated_percentage [Float]"
$ws1.Range("T1").Value  = "sales_team.2.incentives [Currency]"
$ws1.Range("U1").Value  = "sales_team.3.allocated_percentage [Float]"
$ws1.Range("V1").Value  = "sales_team.3.incentives [Currency]"
$ws1.Range("W1").Value  = "sales_team.4.allocated_percentage [Float]"
$ws1.Range("X1").Value  = "sales_team.4.incentives [Currency]"
$ws1.Range("Y1").Value  = "sales_team.5.allocated_percentage [Float]"
$ws1.Range("Z1").Value  = "sales_team.5.incentives [Currency]"

$ws1.Range("AA1").Value = "default_commission_rate [Float]"
$ws1.Range("AB1").Value = "so_required [Check]"
$ws1.Range("AC1").Value = "dn_required [Check]"
$ws1.Range("AD1").Value = "is_frozen [Check]"
$ws1.Range("AE1").Value = "disabled [Check]"

# --- Add the Instructions sheet right after Template --------------------

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Instructions"

# Copy the same bold header style onto the title cell.
$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("A1").Value = "Instructions"
$ws2.Range("A2").Value = "For child tables:"
$ws2.Range("A3").Value = "- Each child table has 5 sets of columns"
$ws2.Range("A4").Value = "- Column format: tablename.row_number.fieldname"
$ws2.Range("A5").Value = "- Example: items.1.item_name, items.2.item_name"
$ws2.Range("A6").Value = "- Leave cells empty if not needed"

$ws1.Activate()
$ws1.Range("A1").Select()
